# Scheduled market-data refresh: update Leve profit calc columns (H-N) per sheet/row.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3465508.2
$ws.Range("I74").Value = 4329852
$ws.Range("J74").Value = 8133.3335
$ws.Range("K74").Value = 4329852
$ws.Range("L74").Value = 8133.3335
$ws.Range("M74").Value = -4328916
$ws.Range("N74").Value = -10005.3335

$ws.Range("H77").Value = 3465508.2
$ws.Range("I77").Value = 4329852
$ws.Range("J77").Value = 8133.3335
$ws.Range("K77").Value = 21649260
$ws.Range("L77").Value = 40666.6675
$ws.Range("M77").Value = -21644580
$ws.Range("N77").Value = -50026.6675

$ws.Range("H129").Value = 1482508.1
$ws.Range("J129").Value = 1950561.1
$ws.Range("L129").Value = 5851683.300000001
$ws.Range("N129").Value = -5861683.300000001

$ws.Range("H137").Value = 1009.4545
$ws.Range("I137").Value = 721.7143
$ws.Range("J137").Value = 1513
$ws.Range("K137").Value = 2165.1429
$ws.Range("L137").Value = 4539
$ws.Range("M137").Value = 384.8571000000002
$ws.Range("N137").Value = -9639

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2370452
$ws.Range("I32").Value = 6498.3022
$ws.Range("K32").Value = 6498.3022
$ws.Range("M32").Value = -6211.3022

$ws.Range("H61").Value = 1588.6666
$ws.Range("I61").Value = 1516.25
$ws.Range("J61").Value = 1694
$ws.Range("K61").Value = 1516.25
$ws.Range("L61").Value = 1694
$ws.Range("M61").Value = -1304.25
$ws.Range("N61").Value = -2118

$ws.Range("H63").Value = 1997
$ws.Range("I63").Value = 1996.7273
$ws.Range("J63").Value = 1998
$ws.Range("K63").Value = 1996.7273
$ws.Range("L63").Value = 1998
$ws.Range("M63").Value = -1310.7273
$ws.Range("N63").Value = -3370

$ws.Range("H66").Value = 1997
$ws.Range("I66").Value = 1996.7273
$ws.Range("J66").Value = 1998
$ws.Range("K66").Value = 9983.636500000001
$ws.Range("L66").Value = 9990
$ws.Range("M66").Value = -6551.636500000001
$ws.Range("N66").Value = -16854

$ws.Range("H74").Value = 1051.7826
$ws.Range("I74").Value = 993.58826
$ws.Range("J74").Value = 1216.6666
$ws.Range("K74").Value = 993.58826
$ws.Range("L74").Value = 1216.6666
$ws.Range("M74").Value = -119.58826
$ws.Range("N74").Value = -2964.6666

$ws.Range("H77").Value = 1051.7826
$ws.Range("I77").Value = 993.58826
$ws.Range("J77").Value = 1216.6666
$ws.Range("K77").Value = 4967.9413
$ws.Range("L77").Value = 6083.333000000001
$ws.Range("M77").Value = -599.9413000000004
$ws.Range("N77").Value = -14819.333

$ws.Range("H122").Value = 1286.6666
$ws.Range("I122").Value = 928.375
$ws.Range("K122").Value = 2785.125
$ws.Range("M122").Value = -335.125

$ws.Range("H136").Value = 1588.6666
$ws.Range("I136").Value = 1516.25
$ws.Range("J136").Value = 1694
$ws.Range("K136").Value = 4548.75
$ws.Range("L136").Value = 5082
$ws.Range("M136").Value = -1998.75
$ws.Range("N136").Value = -10182

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 60000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 60000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 60000
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -60982

$ws.Range("H105").Value = 4381.3438
$ws.Range("I105").Value = 3976.2307
$ws.Range("K105").Value = 3976.2307
$ws.Range("M105").Value = -2229.2307

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2733
$ws.Range("I132").Value = 1505.5714
$ws.Range("J132").Value = 3807
$ws.Range("K132").Value = 4516.7142
$ws.Range("L132").Value = 11421
$ws.Range("M132").Value = -1986.7142
$ws.Range("N132").Value = -16481

$ws.Range("H140").Value = 51803.637
$ws.Range("J140").Value = 51803.637
$ws.Range("L140").Value = 51803.637
$ws.Range("N140").Value = -62163.637

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 301.55554
$ws.Range("I7").Value = 190.66667
$ws.Range("J7").Value = 523.3333
$ws.Range("K7").Value = 572.00001
$ws.Range("L7").Value = 1569.9999
$ws.Range("M7").Value = -460.00001
$ws.Range("N7").Value = -1793.9999

$ws.Range("H25").Value = 5000
$ws.Range("J25").Value = 5000
$ws.Range("L25").Value = 15000
$ws.Range("N25").Value = -15338

$ws.Range("H30").Value = 5000
$ws.Range("J30").Value = 5000
$ws.Range("L30").Value = 15000
$ws.Range("N30").Value = -15204

$ws.Range("H80").Value = 8375
$ws.Range("J80").Value = 10000
$ws.Range("L80").Value = 30000
$ws.Range("N80").Value = -31872

$ws.Range("H83").Value = 8375
$ws.Range("J83").Value = 10000
$ws.Range("L83").Value = 90000
$ws.Range("N83").Value = -99360

$ws.Range("H122").Value = 527233.3
$ws.Range("J122").Value = 1001224.9
$ws.Range("L122").Value = 9011024.1
$ws.Range("N122").Value = -9015924.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 600331.9399999999
$ws.Range("I122").Value = 878626.75
$ws.Range("J122").Value = 3986
$ws.Range("K122").Value = 2635880.25
$ws.Range("L122").Value = 11958
$ws.Range("M122").Value = -2633430.25
$ws.Range("N122").Value = -16858

$ws.Range("H126").Value = 3872.2222
$ws.Range("I126").Value = 4275
$ws.Range("J126").Value = 3066.6667
$ws.Range("K126").Value = 12825
$ws.Range("L126").Value = 9200.000100000001
$ws.Range("M126").Value = -10355
$ws.Range("N126").Value = -14140.0001

$ws.Range("H132").Value = 5380.4443
$ws.Range("I132").Value = 11006
$ws.Range("J132").Value = 3773.1428
$ws.Range("K132").Value = 33018
$ws.Range("L132").Value = 11319.4284
$ws.Range("M132").Value = -30488
$ws.Range("N132").Value = -16379.4284

$ws.Range("H139").Value = 43239.5
$ws.Range("J139").Value = 43239.5
$ws.Range("L139").Value = 43239.5
$ws.Range("N139").Value = -53519.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1124842.8
$ws.Range("I40").Value = 1445412.1
$ws.Range("J40").Value = 2850
$ws.Range("K40").Value = 1445412.1
$ws.Range("L40").Value = 2850
$ws.Range("M40").Value = -1445276.1
$ws.Range("N40").Value = -3122

$ws.Range("H93").Value = 1373.1154
$ws.Range("I93").Value = 817.1667
$ws.Range("J93").Value = 2624
$ws.Range("K93").Value = 817.1667
$ws.Range("L93").Value = 2624
$ws.Range("M93").Value = 430.8333
$ws.Range("N93").Value = -5120

$ws.Range("H100").Value = 5557083.5
$ws.Range("I100").Value = 6945966.5
$ws.Range("J100").Value = 1550.5
$ws.Range("K100").Value = 6945966.5
$ws.Range("L100").Value = 1550.5
$ws.Range("M100").Value = -6945425.5
$ws.Range("N100").Value = -2632.5
